$d = $word.ActiveDocument
function TestPos($pos) {
    $rr = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TB_$pos", $rr)
}
TestPos 50
TestPos 100
TestPos 200
TestPos 315
TestPos 400
Write-Output "done"
